$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "58.675.70"
Set-TextValue "E2" "  +1.49%  "
Set-TextValue "D3" "3.153.11"
Set-TextValue "E3" "  +1.28%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "532.29"
Set-TextValue "E5" "  +0.59%  "
Set-TextValue "D6" "139.92"
Set-TextValue "E6" "  +1.43%  "
Set-TextValue "E7" "  +0.15%  "
Set-TextValue "D8" "0.521"
Set-TextValue "E8" "  +12.54%  "
Set-TextValue "D9" "7.33"
Set-TextValue "E9" "  +0.91%  "
Set-TextValue "E10" "  +2.94%  "
Set-TextValue "E11" "  +4.56%  "
Set-TextValue "E12" "  +2.87%  "
Set-TextValue "D13" "3.695.15"
Set-TextValue "E13" "  +1.14%  "
Set-TextValue "D14" "25.77"
Set-TextValue "E14" "  +1.52%  "
Set-TextValue "E15" "  +5.50%  "
Set-TextValue "D16" "58.693.01"
Set-TextValue "E16" "  +1.46%  "
Set-TextValue "B17" "Polkadot"
Set-TextValue "C17" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D17" "6.22"
Set-TextValue "E17" "  +3.93%  "
Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.144.55"
Set-TextValue "E18" "  +0.78%  "
Set-TextValue "D19" "13.00"
Set-TextValue "E19" "  +3.64%  "
Set-TextValue "E20" "  +2.32%  "
Set-TextValue "D21" "371.48"
Set-TextValue "E21" "  +6.08%  "
Set-TextValue "E22" "  +2.08%  "
Set-TextValue "E23" "  +0.17%  "
Set-TextValue "D24" "69.91"
Set-TextValue "E24" "  +1.37%  "
Set-TextValue "E25" "  +2.86%  "
Set-TextValue "E26" "  +0.14%  "
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.14%  "
Set-TextValue "D28" "8.07"
Set-TextValue "E28" "  +12.36%  "
Set-TextValue "D29" "0.0₃0861"
Set-TextValue "E29" "  -0.57%  "
Set-TextValue "E30" "  +1.01%  "
Set-TextValue "D31" "6.08"
Set-TextValue "D32" "21.83"
Set-TextValue "E32" "  +3.10%  "
Set-TextValue "D33" "5.19"
Set-TextValue "E33" "  +5.76%  "
Set-TextValue "E34" "  +2.10%  "
Set-TextValue "D35" "159.18"
Set-TextValue "E35" "  +0.21%  "
Set-TextValue "E36" "  +4.09%  "
Set-TextValue "E37" "  +8.07%  "
Set-TextValue "D38" "25.26"
Set-TextValue "E38" "  -1.80%  "
Set-TextValue "D39" "2.666.29"
Set-TextValue "E39" "  +11.68%  "
Set-TextValue "E40" "  +1.84%  "
Set-TextValue "D41" "0.0682"
Set-TextValue "E41" "  +2.29%  "
Set-TextValue "E42" "  +5.38%  "
Set-TextValue "E43" "  +2.35%  "
Set-TextValue "D44" "38.66"
Set-TextValue "E44" "  +4.13%  "
Set-TextValue "D45" "0.0284"
Set-TextValue "E45" "  +7.11%  "
Set-TextValue "D47" "3.193.34"
Set-TextValue "E47" "  +1.18%  "
Set-TextValue "E48" "  +12.18%  "
Set-TextValue "E49" "  +2.95%  "
Set-TextValue "E50" "  +2.71%  "
Set-TextValue "D51" "20.10"
Set-TextValue "E51" "  +2.29%  "
